{"js": "// Append the \"If you count to one hundred...\" sentence (with a superscript\n// \"th\") to the end of the last paragraph of the document, right after the\n// existing \"...10-110-210-310 etc.\" sentence and before the trailing\n// bookmark.\n\nconst body = context.document.body;\n\n// Locate the end of the existing sentence that the new text should follow.\nconst anchorResults = body.search(\"10-110-210-310 etc.\", { matchCase: true });\nanchorResults.load(\"items\");\nawait context.sync();\nconst anchor = anchorResults.items[0];\n\n// Insert the first new sentence fragment right after the anchor text.\nconst run1 = anchor.insertText(\n  \" If you count to one hundred you always hit the 100\",\n  \"After\"\n);\nawait context.sync();\n\n// Insert the closing fragment right after run1 (the superscript \"th\" will be\n// inserted between the two in the next step).\nconst run3 = run1.insertText(\n  \" number on your ring finger: whether, it be 100-200-500-1,000 etc.\",\n  \"After\"\n);\nawait context.sync();\n\n// Insert the superscript \"th\" between run1 and run3.\nconst run2 = run1.insertText(\"th\", \"After\");\nrun2.font.superscript = true;\nawait context.sync();\n\n// Force the boundary between the original sentence and the newly inserted\n// text to remain a distinct run (instead of being silently re-merged with\n// the pre-existing run because the two share identical formatting) by\n// toggling a formatting property on run1 and then reverting it. This must\n// happen last, after every insertion above, since further insertText calls\n// on neighbouring ranges otherwise re-normalize (merge) adjacent runs that\n// share identical formatting.\nrun1.font.bold = true;\nawait context.sync();\nrun1.font.bold = false;\nawait context.sync();\n", "ps1": "# Append the \"If you count to one hundred...\" sentence (with a superscript\n# \"th\") to the end of the last paragraph, right after the existing\n# \"...10-110-210-310 etc.\" sentence and before the trailing bookmark.\n\n$d = $word.ActiveDocument\n\n# Locate the end of the existing sentence that the new text should follow.\n$rng = $d.Content\n$rng.Find.Execute(\"10-110-210-310 etc.\")\n$rng.Collapse(0)   # wdCollapseEnd\n$run1Start = $rng.Start\n\n# Insert the first new sentence fragment right after the anchor text. Word\n# merges this into the pre-existing run (same formatting) for now; the\n# boundary is re-established at the very end of the script.\n$rng.Text = \" If you count to one hundred you always hit the 100\"\n$run1End = $rng.End\n\n# Insert the closing (non-superscript) fragment immediately after run1's\n# text *before* the superscript \"th\" is inserted, so it never inherits the\n# superscript formatting.\n$run3Range = $d.Range($run1End, $run1End)\n$run3Range.Text = \" number on your ring finger: whether, it be 100-200-500-1,000 etc.\"\n\n# Insert the superscript \"th\" between run1 and run3.\n$run2Range = $d.Range($run1End, $run1End)\n$run2Range.Text = \"th\"\n$run2Range.Font.Superscript = $true\n\n# Force the boundary between the original sentence and the newly inserted\n# text to remain a distinct run (instead of being silently re-merged with\n# the pre-existing run because the two share identical formatting) by\n# toggling a formatting property on run1's range and then reverting it.\n# This must happen last, after every insertion above.\n$run1Range = $d.Range($run1Start, $run1End)\n$run1Range.Font.Bold = $true\n$run1Range.Font.Bold = $false\n"}
